$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = 'Normal'
}

$ws.Range('D2').Value = '67.599.63'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '3.333.49'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('E4').Value = '  +0.11%  '
Set-TextValue $ws.Range('D5') '580.47'
$ws.Range('E5').Value = '  -1.04%  '
Set-TextValue $ws.Range('D6') '175.83'
$ws.Range('E6').Value = '  -3.58%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '3.330.10'
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  -0.65%  '
Set-TextValue $ws.Range('D12') '45.46'
$ws.Range('E12').Value = '  -2.30%  '
Set-TextValue $ws.Range('D13') '0.0000270'
$ws.Range('E13').Value = '  -2.41%  '
Set-TextValue $ws.Range('D14') '665.57'
$ws.Range('E14').Value = '  +3.41%  '
$ws.Range('D15').Value = '3.879.25'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '67.744.17'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').Value = '3.334.79'
$ws.Range('E19').Value = '  +0.33%  '
Set-TextValue $ws.Range('D20') '17.38'
$ws.Range('E20').Value = '  -1.81%  '
Set-TextValue $ws.Range('D21') '10.99'
$ws.Range('E21').Value = '  +0.28%  '
Set-TextValue $ws.Range('D22') '0.889'
$ws.Range('E22').Value = '  -1.25%  '
Set-TextValue $ws.Range('D23') '5.44'
$ws.Range('E23').Value = '  +8.64%  '
$ws.Range('E24').Value = '  -3.74%  '
Set-TextValue $ws.Range('D25') '99.34'
$ws.Range('E25').Value = '  +1.27%  '
Set-TextValue $ws.Range('D26') '3.86'
$ws.Range('E26').Value = '  -3.59%  '
Set-TextValue $ws.Range('D27') '2.66'
$ws.Range('E27').Value = '  -4.82%  '
Set-TextValue $ws.Range('D28') '9.30'
$ws.Range('E28').Value = '  -3.36%  '
Set-TextValue $ws.Range('D29') '33.63'
$ws.Range('E29').Value = '  +2.09%  '
Set-TextValue $ws.Range('D30') '7.40'
$ws.Range('E30').Value = '  +11.01%  '
$ws.Range('E31').Value = '  -1.78%  '
Set-TextValue $ws.Range('D32') '578.64'
$ws.Range('E32').Value = '  -4.53%  '
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '3.698.01'
$ws.Range('E36').Value = '  -5.52%  '
Set-TextValue $ws.Range('D37') '56.61'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('E38').Value = '  -6.19%  '
Set-TextValue $ws.Range('D39') '34.38'
$ws.Range('E39').Value = '  +4.34%  '
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('E41').Value = '  -2.92%  '
$ws.Range('E42').Value = '  -5.05%  '
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').Value = '0.0₃0667'
$ws.Range('E44').Value = '  -3.11%  '
$ws.Range('E45').Value = '  -2.41%  '
$ws.Range('E46').Value = '  -2.51%  '
Set-TextValue $ws.Range('D47') '2.60'
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('E50').Value = '  +0.51%  '
Set-TextValue $ws.Range('D51') '128.92'
$ws.Range('E51').Value = '  -1.70%  '
